# Weekly update: a new price observation (most recent date) is inserted at
# the top of the data block (row 16, just below the header rows), which
# pushes every existing observation down by one row. The oldest observation
# that used to sit at the bottom of the sheet is preserved (it simply moves
# from row 88 to row 89).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 16 - shifts rows 16:88 down to 17:89 and
# grows the sheet's used range to A1:R89.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row with the latest observation.
$ws.Range("A16").Value = 4
$ws.Range("B16").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C16").Value = "Los Lagos"
$ws.Range("D16").Value = 45030
$ws.Range("E16").Value = 10
$ws.Range("F16").Value = 100112043
$ws.Range("G16").Value = "Pepino dulce"
$ws.Range("H16").Value = "Cultivar IV Región"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = 18000
$ws.Range("L16").Value = 18000
$ws.Range("M16").Value = 18000
$ws.Range("N16").Value = "`$/bandeja 18 kilos"
$ws.Range("O16").Value = "Provincia de Limarí"
$ws.Range("P16").Value = 1000
$ws.Range("Q16").Value = 18
$ws.Range("R16").Value = "Hortaliza"
